# Commit: "update Alfarouk on Adham"
#
# The single paragraph currently reads "svsvs" (plus a leftover
# "_GoBack" bookmark). It needs to become three lines of "typed" text:
#   Adham
#   Alfarouk Updated
# ("Alfarouk" is additionally wrapped in spell-check proofErr markers,
# as Word does for a word it doesn't recognize, and the line break
# between "Adham" and "Alfarouk" is a manual line break, not a new
# paragraph.) The stray "_GoBack" bookmark from the old edit session
# is dropped along the way.

$d = $word.ActiveDocument

$newParagraphXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Adham</w:t></w:r><w:r><w:br/></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Alfarouk</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Updated </w:t></w:r></w:p>
'@

# Replacing the whole story's content with explicit OOXML gives us exact
# control over run/break/proofErr boundaries (and drops the old
# bookmark along with the replaced content).
[void]$d.Content.InsertXML($newParagraphXml)
